$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "/portal/legalAct/lt/TAK/7561b82110ba11e88a05839ea3846d8e?jfwid=-fa58i7sxd"
$ws.Range("A2").Value = "/portal/legalAct/lt/TAK/a3497810124311e88a05839ea3846d8e?jfwid=-fa58i7swi"
